$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark that currently sits alone in an
#    empty paragraph (it will be re-created further below, in the
#    middle of the HTML/CSS/Bootstrap paragraph).
# ------------------------------------------------------------------
$bookmarks = $d.Bookmarks
if ($bookmarks.Exists("_GoBack")) {
    $bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Locate the paragraph describing the HTML / CSS / Bootstrap stack
#    and rewrite it to also mention jQuery / JS, with the relevant
#    technology names in bold.
# ------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*balises HTML*") {
        $target = $para
        break
    }
}

# --- 2a. Insert the new "et un peux de jQuery /et JS" sentence ------
$r = $target.Range
$r.Find.Execute("e balises HTML . pour", $true, $false, $false, $false, $false, $true, 1, $false, "e balises HTML et un peux de jQuery /et JS. pour", 2) | Out-Null

# --- 2b. Bold "jQuery" -----------------------------------------------
$r = $target.Range
$r.Find.Execute("jQuery", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Font.Bold = 1

# --- 2c. Bold "JS" (the one that follows " /et ") --------------------
$r = $target.Range
$r.Find.Execute(" /et JS.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$jsStart = $r.Text.IndexOf("JS")
$jsRange = $d.Range($r.Start + $jsStart, $r.Start + $jsStart + 2)
$jsRange.Font.Bold = 1

# --- 2d. Bold "CSS" ----------------------------------------------------
$r = $target.Range
$r.Find.Execute("CSS", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Font.Bold = 1

# --- 2e. Bold "Bootstarp" (split as Boo + t + starp in the XML, but
#         visually it is simply the single bold word "Bootstarp") ------
$r = $target.Range
$r.Find.Execute("Bootstarp", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Font.Bold = 1

# --- 2f. Drop the trailing period directly glued to "Bootstarp" and put
#         back ". " afterwards (matches the diff: "starp" + ". ") -------
$r = $target.Range
$r.Find.Execute("Bootstarp.", $true, $false, $false, $false, $false, $true, 1, $false, "Bootstarp", 2) | Out-Null

$r = $target.Range
$r.Find.Execute("Bootstarp", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterBold = $d.Range($r.End, $r.End + 2)
if ($afterBold.Text -ne ". ") {
    $afterBold = $d.Range($r.End, $r.End)
    $afterBold.InsertAfter(". ")
}
$afterBold.Font.Bold = 0

# ------------------------------------------------------------------
# 3) Re-insert the _GoBack bookmark right after "...structure et le"
#    (before " design de la plateforme ..."), matching the diff.
# ------------------------------------------------------------------
$r = $target.Range
$r.Find.Execute("structure et le", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bookmarkPoint = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $bookmarkPoint) | Out-Null
